$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting existing row 2 down to row 3.
$ws.Rows.Item(2).Insert()

# New row 2: A2 = 3 (same style as the row that got pushed down, i.e. same as A3), B2 = 2
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 2

# Row 3 (previously row 2, shifted down by Insert) keeps A3 = 0 but B3 changes from 3 to 1
$ws.Range("B3").Value = 1

# Make sure A2 carries the same style (bordered/bold/centered) as A3
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Re-apply the value since PasteSpecial formats shouldn't touch it, but ensure consistency
$ws.Range("A2").Value = 3

# The row insert caused B2 to inherit formatting (bold/centered) from B1 above it;
# B2 should remain plain/unstyled like the original B2 cell, so clear its formatting.
$ws.Range("B2").ClearFormats()
$ws.Range("B2").Value = 2
